# Fixed import function for new equipment categories
#
# The previous import run added a new equipment-category row (ZZZZ / Test)
# without the proper default row height, and this run adds the next row the
# same importer produces (YYYY / Tewst). Re-normalize the older rows' height
# the same way the importer's newly-added rows already look (no explicit
# row height override) and leave the selection where the importer left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newly-imported equipment category as row 28.
$ws.Range("A28").Value = "YYYY"
$ws.Range("B28").Value = "Tewst"

# Clear the stale explicit row heights on the pre-existing rows so they
# follow the sheet's default row height again (matches how freshly
# imported rows, e.g. row 27/28, already behave).
$ws.Rows("1:26").AutoFit()

# Leave the selection where the import routine would land after appending
# the new row.
$ws.Range("E29").Select()
